$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]0x0B

# Cell (1,1): '96 x 68' -> '78 x 45'
$cell = $t.Cell(1,1)
$cell.Range.Text = "78 x 45" + $vt + "  4    5" + $vt + "  ----" + $vt + "7|    |" + $vt + "8|    |"

# Cell (1,2): '58 x 82' -> '62 x 11'
$cell = $t.Cell(1,2)
$cell.Range.Text = "62 x 11" + $vt + "  1    1" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"

# Cell (1,3): '53 x 21' -> '53 x 79'
$cell = $t.Cell(1,3)
$cell.Range.Text = "53 x 79" + $vt + "  7    9" + $vt + "  ----" + $vt + "5|    |" + $vt + "3|    |"

# Cell (2,1): '53 x 41' -> '22 x 13'
$cell = $t.Cell(2,1)
$cell.Range.Text = "22 x 13" + $vt + "  1    3" + $vt + "  ----" + $vt + "2|    |" + $vt + "2|    |"

# Cell (2,2): '89 x 55' -> '11 x 90'
$cell = $t.Cell(2,2)
$cell.Range.Text = "11 x 90" + $vt + "  9    0" + $vt + "  ----" + $vt + "1|    |" + $vt + "1|    |"

# Cell (2,3): '72 x 48' -> '72 x 97'
$cell = $t.Cell(2,3)
$cell.Range.Text = "72 x 97" + $vt + "  9    7" + $vt + "  ----" + $vt + "7|    |" + $vt + "2|    |"

# Cell (3,1): '18 x 30' -> '78 x 24'
$cell = $t.Cell(3,1)
$cell.Range.Text = "78 x 24" + $vt + "  2    4" + $vt + "  ----" + $vt + "7|    |" + $vt + "8|    |"

# Cell (3,2): '83 x 68' -> '80 x 96'
$cell = $t.Cell(3,2)
$cell.Range.Text = "80 x 96" + $vt + "  9    6" + $vt + "  ----" + $vt + "8|    |" + $vt + "0|    |"

# Cell (3,3): '22 x 62' -> '38 x 43'
$cell = $t.Cell(3,3)
$cell.Range.Text = "38 x 43" + $vt + "  4    3" + $vt + "  ----" + $vt + "3|    |" + $vt + "8|    |"

# Cell (4,1): '47 x 76' -> '89 x 60'
$cell = $t.Cell(4,1)
$cell.Range.Text = "89 x 60" + $vt + "  6    0" + $vt + "  ----" + $vt + "8|    |" + $vt + "9|    |"

# Cell (4,2): '33 x 21' -> '61 x 34'
$cell = $t.Cell(4,2)
$cell.Range.Text = "61 x 34" + $vt + "  3    4" + $vt + "  ----" + $vt + "6|    |" + $vt + "1|    |"

# Cell (4,3): '78 x 90' -> '35 x 40'
$cell = $t.Cell(4,3)
$cell.Range.Text = "35 x 40" + $vt + "  4    0" + $vt + "  ----" + $vt + "3|    |" + $vt + "5|    |"

# Cell (5,1): '18 x 37' -> '35 x 49'
$cell = $t.Cell(5,1)
$cell.Range.Text = "35 x 49" + $vt + "  4    9" + $vt + "  ----" + $vt + "3|    |" + $vt + "5|    |"

# Cell (5,2): '75 x 85' -> '92 x 78'
$cell = $t.Cell(5,2)
$cell.Range.Text = "92 x 78" + $vt + "  7    8" + $vt + "  ----" + $vt + "9|    |" + $vt + "2|    |"

# Cell (5,3): '63 x 12' -> '56 x 34'
$cell = $t.Cell(5,3)
$cell.Range.Text = "56 x 34" + $vt + "  3    4" + $vt + "  ----" + $vt + "5|    |" + $vt + "6|    |"

